# ===========================================================
# Edit: add the "2022-Q4" worksheet with the quarter's fund
# holding data, and update the "总计" summary sheet so its
# quarter-over-quarter totals table includes the new quarter.
# ===========================================================

$wb = $excel.ActiveWorkbook

function Set-IndexStyle($cell) {
    # Mirrors the bold / centered / thin-bordered look already used for
    # the header row and the leading index column on every quarterly
    # sheet (and on the "总计" sheet's index column) in this workbook.
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# -----------------------------------------------------------
# 1) Update the "总计" summary sheet (sheet index 1): shift the
#    existing quarterly rows down by one and insert the new
#    2022-Q4 figures at the top of the data table.
# -----------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$totalRows = @(
    @(0, "2022-Q4", 27, 7.28),
    @(1, "2022-Q3", 26, 6.17),
    @(2, "2022-Q2", 6,  0.77),
    @(3, "2022-Q1", 5,  0.84),
    @(4, "2021-Q3", 4,  0.22),
    @(5, "2021-Q2", 2,  0.3),
    @(6, "2021-Q1", 2,  0.06),
    @(7, "2020-Q4", 2,  0.06)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $aCell = $wsTotal.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Set-IndexStyle $aCell

    $wsTotal.Cells.Item($r, 2).Value = $row[1]
    $wsTotal.Cells.Item($r, 3).Value = $row[2]
    $wsTotal.Cells.Item($r, 4).Value = $row[3]
}

# -----------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right after "总计"
#    (i.e. before the current "2022-Q3" sheet) and populate it
#    with the fund-holding breakdown for the quarter.
# -----------------------------------------------------------
$wsQ3Old = $wb.Worksheets.Item(2)   # currently "2022-Q3"; new sheet goes before it
$wsQ4 = $wb.Worksheets.Add($wsQ3Old, $null)
$wsQ4.Name = "2022-Q4"

# Header row
$q4Headers = @(
    "基金代码",
    "基金名称",
    "基金规模",
    "股票总仓位",
    "仓位占比",
    "持有市值(亿元)",
    "仓位排名"
)
for ($col = 2; $col -le 8; $col++) {
    $c = $wsQ4.Cells.Item(1, $col)
    $c.Value = $q4Headers[$col - 2]
    Set-IndexStyle $c
}

# Data rows: columns are
#   A index(0-based,int) | B code(text) | C name(text) | D scale(text)
#   E stock-position(text) | F position-pct(text) | G held-value | H rank(int)
# G is stored as text for every row except the very last one, which the
# source workbook records as the literal number 0.
$q4Data = @(
    @(0, "010190", "嘉实价值发现三个月定期开放混合", "38.99", "95.31", "5.16", "2.0119", 6),
    @(1, "011518", "嘉实价值臻选混合", "27.37", "92.94", "5.75", "1.5738", 8),
    @(2, "070019", "嘉实价值优势混合A", "22.22", "94.25", "5.80", "1.2888", 6),
    @(3, "160605", "鹏华中国50混合", "13.70", "82.57", "3.67", "0.5028", 9),
    @(4, "009989", "华宝研究精选混合", "7.42", "87.41", "4.50", "0.3339", 1),
    @(5, "012262", "华宝可持续发展混合A", "8.23", "89.78", "3.94", "0.3243", 2),
    @(6, "005313", "万家中证1000指数增强A", "22.07", "94.13", "1.11", "0.2450", 1),
    @(7, "005314", "万家中证1000指数增强C", "19.61", "94.13", "1.11", "0.2177", 1),
    @(8, "013624", "嘉实价值创造三年持有期混合A", "3.23", "93.70", "6.31", "0.2038", 5),
    @(9, "012263", "华宝可持续发展混合C", "4.38", "89.78", "3.94", "0.1726", 2),
    @(10, "000866", "华宝高端制造股票", "3.48", "86.44", "3.06", "0.1065", 2),
    @(11, "004845", "南华瑞盈混合A", "2.57", "93.71", "3.00", "0.0771", 10),
    @(12, "013625", "嘉实价值创造三年持有期混合C", "0.75", "93.70", "6.31", "0.0473", 5),
    @(13, "160645", "鹏华精选回报三年定期开放混合", "0.79", "74.91", "4.64", "0.0367", 4),
    @(14, "163110", "申万菱信量化小盘股票（LOF）A", "5.10", "92.35", "0.55", "0.0280", 10),
    @(15, "014839", "兴银碳中和主题混合C", "0.64", "92.17", "3.62", "0.0232", 9),
    @(16, "004284", "华宝新优选一年定期开放灵活配置混合", "0.52", "89.00", "3.69", "0.0192", 3),
    @(17, "014838", "兴银碳中和主题混合A", "0.53", "92.17", "3.62", "0.0192", 9),
    @(18, "002137", "诺安利鑫灵活配置混合A", "0.44", "89.87", "3.80", "0.0167", 6),
    @(19, "002456", "招商安元灵活配置混合A", "0.63", "38.01", "1.62", "0.0102", 8),
    @(20, "002457", "招商安元灵活配置混合C", "0.48", "38.01", "1.62", "0.0078", 8),
    @(21, "015466", "太平中证1000指数增强A", "0.62", "93.58", "0.91", "0.0056", 8),
    @(22, "016169", "嘉实价值优势混合C", "0.09", "94.25", "5.80", "0.0052", 6),
    @(23, "004846", "南华瑞盈混合C", "0.07", "93.71", "3.00", "0.0021", 10),
    @(24, "015467", "太平中证1000指数增强C", "0.08", "93.58", "0.91", "0.0007", 8),
    @(25, "014521", "诺安利鑫灵活配置混合C", "0.01", "89.87", "3.80", "0.0004", 6),
    @(26, "013918", "申万菱信量化小盘股票（LOF）C", "0.00", "92.35", "0.55", 0, 10)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]

    $aCell = $wsQ4.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Set-IndexStyle $aCell

    $wsQ4.Cells.Item($r, 2).NumberFormat = "@"
    $wsQ4.Cells.Item($r, 2).Value = $row[1]

    $wsQ4.Cells.Item($r, 3).Value = $row[2]

    $wsQ4.Cells.Item($r, 4).NumberFormat = "@"
    $wsQ4.Cells.Item($r, 4).Value = $row[3]

    $wsQ4.Cells.Item($r, 5).NumberFormat = "@"
    $wsQ4.Cells.Item($r, 5).Value = $row[4]

    $wsQ4.Cells.Item($r, 6).NumberFormat = "@"
    $wsQ4.Cells.Item($r, 6).Value = $row[5]

    $gVal = $row[6]
    if ($gVal -is [string]) {
        $wsQ4.Cells.Item($r, 7).NumberFormat = "@"
        $wsQ4.Cells.Item($r, 7).Value = $gVal
    } else {
        $wsQ4.Cells.Item($r, 7).Value = $gVal
    }

    $wsQ4.Cells.Item($r, 8).Value = $row[7]
}

